$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 gains the same "data row" formatting used by the rows above it
#     (centered alignment + fixed row height), matching rows 4-13.
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 12.75

# --- New data rows 15-18: same style/height as the existing "data" rows.
$ws.Range("A15").Value = 0.7309
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "20/11/2025"
$ws.Range("D15").Value = "19:00:23"

$ws.Range("A16").Value = 0.7214
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = "20/11/2025"
$ws.Range("D16").Value = "20:14:13"

$ws.Range("A17").Value = 0.7184
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "21/11/2025"
$ws.Range("D17").Value = "01:00:17"

$ws.Range("A18").Value = 0.6979
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "21/11/2025"
$ws.Range("D18").Value = "07:00:15"

# --- New row 19: values only, plain/default formatting (no centering).
$ws.Range("A19").Value = 0.6558
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "21/11/2025"
$ws.Range("D19").Value = "14:00:15"
$ws.Range("G4").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)

# --- New row 20: values only, plain/default formatting (no centering).
$ws.Range("A20").Value = 0.6427
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "22/11/2025"
$ws.Range("D20").Value = "07:01:01"
$ws.Range("G4").Copy()
$ws.Range("A20:D20").PasteSpecial(-4122)

# --- Edited cell: D2 (2 -> 3), then leave the selection there.
$ws.Range("D2").Value = 3
[void]$ws.Range("D2").Select()
